# 2.a.1 - Add a new "2022" column (S) to the agriculture orientation index table
# and refresh the sheet's column widths / active selection to match the
# state the workbook was left in after the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Extend the year header (row 3) and the data row (row 4) with a new
#     column S, copying the formatting of the preceding column R so the new
#     cells line up with the rest of the table. ---
$ws.Range("R3").Copy($ws.Range("S3"))
$ws.Range("S3").Value = 2022

$ws.Range("R4").Copy($ws.Range("S4"))
$ws.Range("S4").Value = 0.071025550219041236

# --- Re-balance the widths of the first three (label) columns so they share
#     a single, equal width instead of three different custom widths. ---
$ws.Range("A:C").ColumnWidth = 32.66

# --- Leave the sheet with the same active cell/selection it had when the
#     file was saved. ---
$ws.Range("F14").Select() | Out-Null
